# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet to reflect newly scraped data.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 334
$wsExpo.Range("F3").Value = 98
$wsExpo.Range("F4").Value = 1383
$wsExpo.Range("F5").Value = 658

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 334
$wsAll.Range("F3").Value = 98
$wsAll.Range("F4").Value = 1383
$wsAll.Range("F6").Value = 658
